# Add a new "2022-Q3" quarter sheet to the workbook, positioned right
# after "总计" and before the existing "2022-Q1" sheet, and add a matching
# summary row on the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new worksheet before the current 2nd sheet ("2022-Q1")
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item(2)        # "2022-Q1" - used as a formatting template
$beforeSheet = $template
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

# Match sheet-level look & feel (outline direction, margins) used by
# every other sheet in this workbook.
$q3.Outline.SummaryRow = 1
$q3.Outline.SummaryColumn = 1
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# Bring over the header-row / index-column formatting (borders, bold,
# centering) by copying the equivalent cells from the template sheet -
# this carries the right style index across without us having to guess it.
$template.Range("A1:H3").Copy($q3.Range("A1:H3"))

# ---- header row text ----
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# ---- row 2 : 005585 银河文体娱乐主题灵活配置混合A ----
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'005585"
$q3.Range("C2").Value = "银河文体娱乐主题灵活配置混合A"
$q3.Range("D2").Value = "'3.01"
$q3.Range("E2").Value = "'90.28"
$q3.Range("F2").Value = "'6.39"
$q3.Range("G2").Value = "'0.1923"
$q3.Range("H2").Value = 2

# ---- row 3 : 015667 银河文体娱乐主题灵活配置混合C ----
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "'015667"
$q3.Range("C3").Value = "银河文体娱乐主题灵活配置混合C"
$q3.Range("D3").Value = "'0.41"
$q3.Range("E3").Value = "'90.28"
$q3.Range("F3").Value = "'6.39"
$q3.Range("G3").Value = "'0.0262"
$q3.Range("H3").Value = 2

# The apostrophe-prefixed assignments above mark those cells as
# "text-with-quote-prefix" - strip that cosmetic marker back off (values
# stay text) by re-pasting the plain formatting from an unstyled cell.
$template.Range("C2").Copy()
$q3.Range("B2:G3").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row for 2022-Q3 right
#    above the existing 2022-Q1 row, pushing the rest down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Shift existing data rows down (old row2->row3, old row3->row4), copying
# so the formatting/style of each source row travels with it.
$total.Range("A3:D3").Copy($total.Range("A4:D4"))
$total.Range("A2:D2").Copy($total.Range("A3:D3"))

# New 2022-Q3 row
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.22

# Fix up the running index column for the rows that shifted down
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
